$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 242; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = "'99999999"
}
